$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36
$ws.Range("A36").Value = "Globo"
$ws.Range("B36").Value = "RJ TV 2"
$ws.Range("C36").Value = "Defesa Civil"
$ws.Range("D36").Value = "2025-04-03T19:29"
$ws.Range("E36").Value = "Positivo"
$ws.Range("F36").Value = "Alerta de chuvas fortes em Campos. Defesa Civil emitiu alerta de fortes chuvas que podem atingir região nos próximos dias. Repórter *ao vivo*. Além de chuvas, baixas temperaturas também. Equipes da Defesa Civil informaram que está sendo esperado acumulado de 100mm a 129mm, considerado volume intenso. Para receber aviso da Defesa Civil, enviar mensagem de texto sem hífen para o 40199."

# Row 37
$ws.Range("A37").Value = "Globo"
$ws.Range("B37").Value = "RJ TV 2"
$ws.Range("C37").Value = "Iluminação"
$ws.Range("D37").Value = "2025-04-03T19:32"
$ws.Range("E37").Value = "Negativo"
$ws.Range("F37").Value = "Por whatsApp, morador de Morro Grande, Lagoa de Cima, reclama da iluminação da rua da localidade. Exibido vídeo enviado por morador. Produção pediu nota à prefeitura e aguarda resposta. *sem nota*"

# Row 38
$ws.Range("A38").Value = "Globo"
$ws.Range("B38").Value = "RJ TV 2"
$ws.Range("C38").Value = "Iluminação"
$ws.Range("D38").Value = "2025-04-03T19:32"
$ws.Range("E38").Value = "Neutro"
$ws.Range("F38").Value = "Poste quebrado há mais de um mês no cruzamento entre a Rua visconde de Itaboraí e a Avenida Arthur Bernardes. Enel enviou resposta. "
